# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (new data for date 2021-11-11,
# serial 44511) above the existing "Frutilla" rows, pushing the rest of
# the table (previously rows 123-155) down to rows 125-157, matching the
# new sheet dimension A1:T157.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 123; everything below shifts down by 2.
$ws.Rows("123:124").Insert()

# Seed the two new rows with the same row "shape" (market/product columns
# A,B,C,E,F,G,H,I,J,K,T) as the row that landed right below them, then
# overwrite the fields that actually change for each new observation.
$ws.Range("A125:T125").Copy()
$ws.Range("A123:T123").PasteSpecial()

$ws.Range("A125:T125").Copy()
$ws.Range("A124:T124").PasteSpecial()

# New row 123: Frutilla, Especial quality
$ws.Cells.Item(123,4).Value = 44511
$ws.Cells.Item(123,12).Value = 'Especial'
$ws.Cells.Item(123,13).Value = 80
$ws.Cells.Item(123,14).Value = 9500
$ws.Cells.Item(123,15).Value = 10000
$ws.Cells.Item(123,16).Value = 9750
$ws.Cells.Item(123,17).Value = '$/caja 7 kilos'
$ws.Cells.Item(123,18).Value = 'Provincia de Diguillín'
$ws.Cells.Item(123,19).Value = 1393

# New row 124: Frutilla, Primera quality
$ws.Cells.Item(124,4).Value = 44511
$ws.Cells.Item(124,12).Value = 'Primera'
$ws.Cells.Item(124,13).Value = 120
$ws.Cells.Item(124,14).Value = 7500
$ws.Cells.Item(124,15).Value = 8000
$ws.Cells.Item(124,16).Value = 7750
$ws.Cells.Item(124,17).Value = '$/caja 7 kilos'
$ws.Cells.Item(124,18).Value = 'Provincia de Diguillín'
$ws.Cells.Item(124,19).Value = 1107
